$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rob Oudman"
$ws.Range("A2").Select()
